$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# The author removed the note "The EU EPS uses values from the US EPS."
# (row 10) together with the blank spacer row beneath it (row 11). Delete
# both rows; everything below shifts up by two rows automatically.
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(10).Delete()

# Reset the lingering cell selection (the source file no longer pins the
# view to C24 once the note above it is gone).
$ws.Range("A1").Select()
